$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.175.99"
$ws.Range("E2").Value = "  +2.07%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.349.65"
$ws.Range("E3").Value = "  +1.26%  "

# Row 4
$ws.Range("E4").Value = "  -0.36%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.26"
$ws.Range("E5").Value = "  +1.85%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.61"
$ws.Range("E6").Value = "  +2.36%  "

# Row 7
$ws.Range("E7").Value = "  +0.59%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.565"
$ws.Range("E8").Value = "  +5.73%  "

# Row 9
$ws.Range("E9").Value = "  +0.20%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.55"

# Row 11
$ws.Range("E11").Value = "  -0.60%  "

# Row 12
$ws.Range("E12").Value = "  +1.86%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.82"
$ws.Range("E13").Value = "  +1.49%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.768.56"
$ws.Range("E14").Value = "  +1.01%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.143.84"
$ws.Range("E15").Value = "  +1.95%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000133"
$ws.Range("E16").Value = "  +0.61%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.363.27"
$ws.Range("E17").Value = "  +1.22%  "

# Row 18
$ws.Range("E18").Value = "  +2.58%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "332.11"
$ws.Range("E19").Value = "  -1.52%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.28"
$ws.Range("E20").Value = "  +2.69%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.83"
$ws.Range("E21").Value = "  -0.89%  "

# Row 22
$ws.Range("E22").Value = "  +0.19%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "62.75"
$ws.Range("E23").Value = "  +1.85%  "

# Row 24
$ws.Range("E24").Value = "  +0.13%  "

# Row 25
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.49"
$ws.Range("E25").Value = "  -2.67%  "

# Row 26
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.51%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.38"
$ws.Range("E27").Value = "  +2.74%  "

# Row 28
$ws.Range("E28").Value = "  +1.63%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.55"
$ws.Range("E29").Value = "  -1.27%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0736"
$ws.Range("E30").Value = "  +1.21%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.14"
$ws.Range("E31").Value = "  +0.46%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.04"
$ws.Range("E32").Value = "  +13.15%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.44"
$ws.Range("E33").Value = "  -0.35%  "

# Row 34
$ws.Range("E34").Value = "  +0.02%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.26"
$ws.Range("E35").Value = "  +6.66%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.71%  "

# Row 37
$ws.Range("E37").Value = "  -0.48%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.64"
$ws.Range("E38").Value = "  +4.72%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.26"
$ws.Range("E39").Value = "  -0.06%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "144.67"
$ws.Range("E40").Value = "  -2.94%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "295.09"
$ws.Range("E41").Value = "  +4.82%  "

# Row 42
$ws.Range("E42").Value = "  +0.54%  "

# Row 43
$ws.Range("E43").Value = "  +1.20%  "

# Row 44
$ws.Range("E44").Value = "  +1.94%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0503"
$ws.Range("E46").Value = "  +0.33%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.562"
$ws.Range("E47").Value = "  +0.37%  "

# Row 48
$ws.Range("E48").Value = "  +1.34%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.51"
$ws.Range("E49").Value = "  +0.24%  "

# Row 50
$ws.Range("E50").Value = "  -0.07%  "

# Row 51
$ws.Range("E51").Value = "  +0.40%  "
